$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.35550000000003
$ws.Range("A21").Value = -21.21220000000001
$ws.Range("A23").Value = -21.34610000000003
$ws.Range("A25").Value = -22.35250000000003
